# Updated cryptos list on Thu Sep  7 21:08:57 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it "looks"
# like a number (e.g. "215.35", "10.00"). Force text number format,
# assign, then restore the default "Normal" style so no stray style
# index is left referenced on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.005.48"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.639.56"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.31%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "215.35"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.504"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.37%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.54%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.0633"
$ws.Range("E9").Value = "  -0.33%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "19.77"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 / Row 13 - Polkadot and WrappedEther swap places
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.668.42"
$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.27"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.862.73"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.554"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  +0.06%  "

# Row 17 - Litecoin
Set-TextValue $ws.Range("D17") "63.08"
$ws.Range("E17").Value = "  +0.68%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.996.08"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19 - Dai
Set-TextValue $ws.Range("D19") "0.999"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "4.46"
$ws.Range("E20").Value = "  +0.42%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "192.28"
$ws.Range("E21").Value = "  -0.77%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("D22") "10.00"
$ws.Range("E22").Value = "  +0.62%  "

# Row 23 - Chainlink
Set-TextValue $ws.Range("D23") "6.37"
$ws.Range("E23").Value = "  +1.60%  "

# Row 24 - BinanceUSD
Set-TextValue $ws.Range("D24") "0.999"
$ws.Range("E24").Value = "  -0.31%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -1.84%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "142.52"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27 - Stellar
Set-TextValue $ws.Range("D27") "0.124"
$ws.Range("E27").Value = "  +1.15%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "6.89"
$ws.Range("E28").Value = "  +0.47%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "15.57"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.27%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.34%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.63%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +0.49%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.68%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.37%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("D36") "0.911"
$ws.Range("E36").Value = "  +1.30%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.148.65"
$ws.Range("E37").Value = "  +2.05%  "

# Row 38 - ImmutableX
Set-TextValue $ws.Range("D38") "0.547"
$ws.Range("E38").Value = "  +0.22%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.08%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.87%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.29%  "

# Row 42 - FraxShare
Set-TextValue $ws.Range("D42") "5.67"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43 - Quant
Set-TextValue $ws.Range("D43") "100.84"

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +0.51%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.773.29"
$ws.Range("E45").Value = "  +0.37%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "55.72"
$ws.Range("E46").Value = "  +1.40%  "

# Row 47 - RenderToken
Set-TextValue $ws.Range("D47") "1.47"
$ws.Range("E47").Value = "  +6.25%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +1.67%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "7.64"
$ws.Range("E50").Value = "  +0.58%  "

# Row 51 - Algorand
Set-TextValue $ws.Range("D51") "0.0970"
$ws.Range("E51").Value = "  +3.36%  "
